$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-234 down to 107-235
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with data (mirrors formatting/values of neighboring rows,
# matching the "Fruta, Feria Lagunitas de Puerto Montt - Mango" weekly data series)
$ws.Range("A106").Value = 4
$ws.Range("B106").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C106").Value = "Los Lagos"
$ws.Range("D106").Value = 44797
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 10
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100108
$ws.Range("H106").Value = "Tropicales y subtropicales"
$ws.Range("I106").Value = 100108002
$ws.Range("J106").Value = "Mango"
$ws.Range("K106").Value = "Sin especificar"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 50
$ws.Range("N106").Value = 13000
$ws.Range("O106").Value = 14000
$ws.Range("P106").Value = 13500
$ws.Range("Q106").Value = "`$/bandeja 4 kilos"
$ws.Range("R106").Value = "Brasil"
$ws.Range("S106").Value = 3375
$ws.Range("T106").Value = 4
